$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows (12-17) below row 11, using row 10 as a style template
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows(10).Copy()
    $ws.Rows(12).Insert(-4121)
}
$excel.CutCopyMode = $false

# Row 12
$ws.Cells.Item(12, 1).Value = 'SCRIPT/T01P02A/um2107.ssb'
$ws.Cells.Item(12, 2).Value = 266
$ws.Cells.Item(12, 3).Value = ' The story is making the rounds\nin Treasure Town. It\''s hard to believe…'
$ws.Cells.Item(12, 4).Value = ' В Городе Сокровищ только и\nговорят об этой истории. В неё трудно\nповерить...'
$ws.Cells.Item(12, 5).Value = ' Â Ãïñïäå Òïëñïâéþ óïìûëï é\nãïâïñÿó ïá üóïê éòóïñéé. Â îåæ óñôäîï\nðïâåñéóû...'

# Row 13
$ws.Cells.Item(13, 1).Value = 'SCRIPT/T01P02A/um2112.ssb'
$ws.Cells.Item(13, 2).Value = 269
$ws.Cells.Item(13, 3).Value = ' But it is true that time hasn\''t\ngone back to normal where stolen Time Gears\nwere put back in place…'
$ws.Cells.Item(13, 4).Value = ' Но действительно, время не\nвернулось в норму когда Шестерни Времени\nвернули на место...'
$ws.Cells.Item(13, 5).Value = ' Îï äåêòóâéóåìûîï, âñåíÿ îå\nâåñîôìïòû â îïñíô ëïãäà Šåòóåñîé Âñåíåîé\nâåñîôìé îà íåòóï...'

# Row 14
$ws.Cells.Item(14, 1).Clear()
$ws.Cells.Item(14, 2).Value = 272
$ws.Cells.Item(14, 3).Value = ' Worse yet, the zones where\ntime has stopped are growing larger…'
$ws.Cells.Item(14, 4).Value = ' И что ещё хуже, зоны, поражённые\nостановкой времени, только растут...'
$ws.Cells.Item(14, 5).Value = ' É œóï åþæ öôçå, èïîú, ðïñàçæîîúå\nïòóàîïâëïê âñåíåîé, óïìûëï ñàòóôó...'

# Row 15
$ws.Cells.Item(15, 1).Clear()
$ws.Cells.Item(15, 2).Value = 275
$ws.Cells.Item(15, 3).Value = ' There are things happening that\ndon\''t make sense…'
$ws.Cells.Item(15, 4).Value = ' Происходит столько всего\nнепонятного...'
$ws.Cells.Item(15, 5).Value = ' Ðñïéòöïäéó òóïìûëï âòåãï\nîåðïîÿóîïãï...'

# Row 16
$ws.Cells.Item(16, 1).ClearContents()
$ws.Cells.Item(16, 2).Value = 278
$ws.Cells.Item(16, 3).Value = ' That\''s why I don\''t think it\''s safe\nto deny the rumor. You know! The rumor\nabout [CS:N]Grovyle[CR]…'
$ws.Cells.Item(16, 4).Value = ' И поэтому я считаю, что слухи\nмогут оказаться правдой. Ну, вы знаете!\nСлухи о [CS:N]Гровайле[CR]...'
$ws.Cells.Item(16, 5).Value = ' É ðïüóïíô ÿ òœéóàý, œóï òìôöé\níïãôó ïëàèàóûòÿ ðñàâäïê. Îô, âú èîàåóå!\nÒìôöé ï [CS:N]Ãñïâàêìå[CR]...'

# Row 17
$ws.Cells.Item(17, 1).Value = 'SCRIPT/T01P02A/um2401.ssb'
$ws.Cells.Item(17, 2).Value = 247
$ws.Cells.Item(17, 3).Value = ' Go with confidence and it\''ll be\nall right! Be tough!'
$ws.Cells.Item(17, 4).Value = ' Идите уверенно и у вас всё\nполучится! Будьте сильными!'
$ws.Cells.Item(17, 5).Value = ' Éäéóå ôâåñåîîï é ô âàò âòæ\nðïìôœéóòÿ! Áôäûóå òéìûîúíé!'

# Fix styles: row 11 and row 16 need the "block separator" style (s=10 / s=11), matching rows 7/9
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A9:E9").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights
$ws.Rows(12).RowHeight = 43.2
$ws.Rows(13).RowHeight = 43.2
$ws.Rows(14).RowHeight = 31.8
$ws.Rows(15).RowHeight = 21.6
$ws.Rows(16).RowHeight = 31.8
$ws.Rows(17).RowHeight = 43.2

# Final selection state
$ws.Range("D17").Select() | Out-Null
